$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (pushes existing rows 6-18 down to 7-19).
# Excel copies formatting from the row above (row 5) automatically.
$ws.Rows.Item(6).Insert()

# New field: summary_malfunction_reporting
$ws.Cells.Item(6, 2).Value = "summary_malfunction_reporting"
$ws.Cells.Item(6, 3).Value = "string"
$ws.Cells.Item(6, 4).Value = "The Voluntary Malfunction Summary Reporting Program allows participating companies to submit certain medical device malfunction reports in summary form on a quarterly basis.  The program applies to eligible devices regulated by the Center for Devices and Radiological Health (CDRH) and Center for Biologics Evaluation and Research (CBER), including device-led combination products.Value is one of the following:
Eligible = 510(K)
Ineligible = PMA"

# Row height for the new row (wrapped text row, matches the other description rows).
$ws.Rows.Item(6).RowHeight = 102

# Column B got a bit wider (and lost its auto "best fit" flag) to fit the new field name.
$ws.Columns.Item(2).ColumnWidth = 27

# Window/view state: zoom in and scroll/select near the new row.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 170
$ws.Range("D6").Select()
